$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.2
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 1.75
$ws.Range("J2").Value = 4.75
$ws.Range("K2").Value = 2.2
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 3.5
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.85
$ws.Range("S2").Value = 1.4
$ws.Range("T2").Value = 2.75
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 1.83
$ws.Range("W2").Value = 11
$ws.Range("AA2").Value = 34
$ws.Range("AB2").Value = 41
$ws.Range("AC2").Value = 11
$ws.Range("AE2").Value = 17
$ws.Range("AF2").Value = 51
$ws.Range("AG2").Value = 301
$ws.Range("AH2").Value = 7
$ws.Range("AI2").Value = 8
$ws.Range("AM2").Value = 26
$ws.Range("AO2").Value = 23
$ws.Range("AP2").Value = 34
$ws.Range("AQ2").Value = 81
$ws.Range("AR2").Value = 101
$ws.Range("AS2").Value = 251
$ws.Range("AT2").Value = 2.75
$ws.Range("AU2").Value = 8.5
$ws.Range("AW2").Value = 3.75
$ws.Range("AY2").Value = 21
$ws.Range("BB2").Value = 151

# Row 3
$ws.Range("Q3").Value = 1.4
$ws.Range("R3").Value = 2.88
$ws.Range("S3").Value = 1.17

# Row 4
$ws.Range("G4").Value = 1.85
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 2.5
$ws.Range("K4").Value = 2.25
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 1.67
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 8.5
$ws.Range("Z4").Value = 17
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 41
$ws.Range("AG4").Value = 201
$ws.Range("AH4").Value = 12
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 12
$ws.Range("AK4").Value = 41
$ws.Range("AL4").Value = 29
$ws.Range("AM4").Value = 34
$ws.Range("AN4").Value = 4
$ws.Range("AP4").Value = 19
$ws.Range("AQ4").Value = 34
$ws.Range("AR4").Value = 51
$ws.Range("AS4").Value = 126
$ws.Range("AT4").Value = 3
$ws.Range("AU4").Value = 7.5
$ws.Range("AV4").Value = 51
$ws.Range("AW4").Value = 5.5
$ws.Range("AX4").Value = 21
$ws.Range("AY4").Value = 26
$ws.Range("AZ4").Value = 67
$ws.Range("BA4").Value = 81
$ws.Range("BB4").Value = 151

# Row 5
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3.2

# Row 8
$ws.Range("N8").Value = 10
$ws.Range("BD8").Value = 126

# Row 10
$ws.Range("G10").Value = 1.7
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 4.2
$ws.Range("J10").Value = 2.38
$ws.Range("L10").Value = 4.5
$ws.Range("N10").Value = 13
$ws.Range("Q10").Value = 1.75
$ws.Range("R10").Value = 2.05
$ws.Range("U10").Value = 1.67
$ws.Range("V10").Value = 2.1
$ws.Range("W10").Value = 8.5
$ws.Range("X10").Value = 9
$ws.Range("Z10").Value = 15
$ws.Range("AE10").Value = 13
$ws.Range("AG10").Value = 151
$ws.Range("AI10").Value = 23
$ws.Range("AO10").Value = 9
$ws.Range("AQ10").Value = 29
$ws.Range("AW10").Value = 6
$ws.Range("AZ10").Value = 67
$ws.Range("BA10").Value = 81
$ws.Range("BB10").Value = 151

# Row 11
$ws.Range("H11").Value = 3.5
$ws.Range("L11").Value = 5.2
$ws.Range("N11").Value = 6.85
$ws.Range("P11").Value = 2.82
$ws.Range("U11").Value = 1.9
$ws.Range("W11").Value = 6
$ws.Range("X11").Value = 7.2
$ws.Range("AE11").Value = 17.5
$ws.Range("AH11").Value = 12.5
$ws.Range("AI11").Value = 29
$ws.Range("AM11").Value = 60
$ws.Range("AN11").Value = 3.4
$ws.Range("AO11").Value = 7.9
$ws.Range("AP11").Value = 17.5
$ws.Range("AQ11").Value = 26
$ws.Range("AR11").Value = 60
$ws.Range("AS11").Value = 250
$ws.Range("AT11").Value = 2.7
$ws.Range("AU11").Value = 7.5
$ws.Range("AV11").Value = 70
$ws.Range("AW11").Value = 6.6
$ws.Range("AY11").Value = 32
$ws.Range("BA11").Value = 200
$ws.Range("BB11").Value = 450
